$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the existing data rows (keeps per-cell styles, drops now-unreferenced
# shared strings) so we can re-populate them in the exact order the target
# shared-strings table needs.
$ws.Range("A2:D11").ClearContents()

# Target rows, in final top-to-bottom order. Re-entering the data in this
# order (row by row, column A,B,C,D) makes every *new* string get interned
# into the shared-strings table in the same sequence the target file uses;
# strings that repeat (e.g. "T1498.001") automatically reuse their earlier slot.
$rows = @(
    @("CyberSec.Booster: ARP spoofing", "T1557.002", "SW_MATM-4 SW_DAI-4", "e9c101d0-f344-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: Password Guessing", "T1110.001", "SSH-5 SEC_LOGIN-4", "d36fef60-f347-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: Manipulation of the STP protocol", "T1498.001", "SPANTREE", "267c1370-f346-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: Detection of a rogue DHCP server", "T1557.003", "SW_DAI-4", "cbfe07a0-f345-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: Scanning ip blocks", "T1595.001", "", "eb158cb0-e54b-11ee-aad9-f582020d7fab"),
    @("CyberSec.Booster: Vulnerability Scanning", "T1595.002", "", "d884bc00-f343-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: DNS data exfiltration", "T1071.004", "", "8adc2400-f345-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: Telnet Port Activity", "TA0011", "", "9610b910-f347-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: Manipulation of the DTP protocol", "T1557", "DTP-5", "0007b730-f346-11ee-95ad-619443a476e1"),
    @("CyberSec.Booster: MAC Flooding", "T1498.001", "MACNOTIFY-6 PORT_SECURITY-2", "0aab2540-f347-11ee-95ad-619443a476e1")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne "") {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# The "facility" style (number-format xf index 2) moves from D9 to D8.
$ws.Range("D9").ClearFormats()
$ws.Range("D8").NumberFormat = "0.00E+00"

# Selection moves off the table, onto C15.
$ws.Range("C15").Select()
